$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.234.41'
$ws.Range('E2').Value = '  -3.77%  '
$ws.Range('D3').Value = '2.244.35'
$ws.Range('E3').Value = '  -4.27%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''245.26'
$ws.Range('E5').Value = '  +2.61%  '
$ws.Range('D6').Value = '''0.631'
$ws.Range('E6').Value = '  -5.27%  '
$ws.Range('D7').Value = '''69.15'
$ws.Range('E7').Value = '  -4.72%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.554'
$ws.Range('E9').Value = '  -6.44%  '
$ws.Range('D10').Value = '''0.0986'
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').Value = '''59.24'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '''36.33'
$ws.Range('E12').Value = '  +11.06%  '
$ws.Range('E13').Value = '  -2.54%  '
$ws.Range('D14').Value = '''6.74'
$ws.Range('E14').Value = '  -7.00%  '
$ws.Range('D15').Value = '2.579.29'
$ws.Range('E15').Value = '  -4.24%  '
$ws.Range('D16').Value = '''15.01'
$ws.Range('E16').Value = '  -6.55%  '
$ws.Range('D17').Value = '''0.865'
$ws.Range('E17').Value = '  -3.88%  '
$ws.Range('D18').Value = '2.247.82'
$ws.Range('E18').Value = '  -3.89%  '
$ws.Range('D19').Value = '42.156.59'
$ws.Range('E19').Value = '  -3.67%  '
$ws.Range('D20').Value = '0.0₃0970'
$ws.Range('E20').Value = '  -6.01%  '
$ws.Range('D21').Value = '''6.24'
$ws.Range('E21').Value = '  -6.18%  '
$ws.Range('D22').Value = '''73.18'
$ws.Range('E22').Value = '  -6.48%  '
$ws.Range('D23').Value = '''235.94'
$ws.Range('E23').Value = '  -6.11%  '
$ws.Range('D24').Value = '''2.04'
$ws.Range('E24').Value = '  +11.31%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''2.23'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''9.99'
$ws.Range('E29').Value = '  -3.79%  '
$ws.Range('D30').Value = '''172.72'
$ws.Range('E30').Value = '  -2.18%  '
$ws.Range('D31').Value = '''20.53'
$ws.Range('E31').Value = '  -7.52%  '
$ws.Range('E32').Value = '  -3.30%  '
$ws.Range('D33').Value = '''0.126'
$ws.Range('E33').Value = '  -5.07%  '
$ws.Range('D34').Value = '''0.0717'
$ws.Range('E34').Value = '  -3.72%  '
$ws.Range('E35').Value = '  -0.57%  '
$ws.Range('E36').Value = '  -7.05%  '
$ws.Range('D37').Value = '''3.79'
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('D38').Value = '''22.81'
$ws.Range('E38').Value = '  +21.68%  '
$ws.Range('D39').Value = '''0.0283'
$ws.Range('E39').Value = '  +4.63%  '
$ws.Range('D40').Value = '''2.30'
$ws.Range('E40').Value = '  -2.81%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '''5.91'
$ws.Range('E41').Value = '  -7.51%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = '''66.96'
$ws.Range('E42').Value = '  +2.57%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').Value = '''5.07'
$ws.Range('E43').Value = '  -9.15%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''9.28'
$ws.Range('E44').Value = '  +0.78%  '
$ws.Range('E45').Value = '  -2.72%  '
$ws.Range('D46').Value = '''4.69'
$ws.Range('E46').Value = '  +12.41%  '
$ws.Range('D47').Value = '''0.190'
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('E49').Value = '  -1.68%  '
$ws.Range('D50').Value = '''10.06'
$ws.Range('E50').Value = '  +7.78%  '
$ws.Range('D51').Value = '''2.33'
$ws.Range('E51').Value = '  -3.14%  '
